# Regenerate merged AHB files
# - Rename header columns: "<name>_old" -> "<name>_FV2210", "<name>_new" -> "<name>_FV2304"
# - Wrap the data range in an Excel Table (ListObject) with AutoFilter, named "Table1"
# - Freeze the header row (pane split at row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newHeaders = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $newHeaders[$i]
}

# Turn the used range into a proper Excel table with an AutoFilter.
$tableRange = $ws.UsedRange
$listObject = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$listObject.Name = "Table1"

# Freeze the header row.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
